$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the stray empty-string cells left over on row 7 (I7, K7:R7)
$ws.Range("I7").ClearContents()
$ws.Range("K7:R7").ClearContents()

# New rows of YSL LOVESHINE PLUMPING GLOSS products (rows 8-17) plus a
# final TESTE row (row 18), added to support the new "bloqueo" (lock)
# review workflow state "Solo Revisión".

$data = @(
    @("0ML34254", "YSL LOVESHINE PLUMPING GLOSS N 4", "MAQUILLAJE LABIOS", "Tiene PT", "Tiene ES", "Tiene IT", "6", "ML", "Solo Revisión"),
    @("0ML34251", "YSL LOVESHINE PLUMPING GLOSS N 1", "MAQUILLAJE LABIOS", "Tiene PT", "Tiene ES", "Tiene IT", "6", "ML", "Solo Revisión"),
    @("0ML34257", "YSL LOVESHINE PLUMPING GLOSS N 8", "MAQUILLAJE LABIOS", "Tiene PT", "Tiene ES", "Tiene IT", "1", "ML", "Solo Revisión"),
    @("0ML34253", "YSL LOVESHINE PLUMPING GLOSS N 3", "MAQUILLAJE LABIOS", "Tiene PT", "Tiene ES", "Tiene IT", "6", "ML", "Solo Revisión"),
    @("0ML34258", "YSL LOVESHINE PLUMPING GLOSS N 44", "MAQUILLAJE LABIOS", "Tiene PT", "Tiene ES", "Tiene IT", "6", "ML", "Solo Revisión"),
    @("0ML34252", "YSL LOVESHINE PLUMPING GLOSS N 2", "MAQUILLAJE LABIOS", "Tiene PT", "Tiene ES", "Tiene IT", "1", "ML", "Solo Revisión"),
    @("0ML34255", "YSL LOVESHINE PLUMPING GLOSS N 6", "MAQUILLAJE LABIOS", "Tiene PT", "Tiene ES", "Tiene IT", "6", "ML", "Solo Revisión"),
    @("0ML34256", "YSL LOVESHINE PLUMPING GLOSS N 7", "MAQUILLAJE LABIOS", "Tiene PT", "Tiene ES", "Tiene IT", "6", "ML", "Solo Revisión"),
    @("0ML34251", "YSL LOVESHINE PLUMPING GLOSS N 1", "MAQUILLAJE LABIOS", "Tiene PT", "Tiene ES", "Tiene IT", "6", "ML", "Solo Revisión"),
    @("0ML34257", "YSL LOVESHINE PLUMPING GLOSS N 8", "MAQUILLAJE LABIOS", "Tiene PT", "Tiene ES", "Tiene IT", "6", "ML", "Solo Revisión")
)

$row = 8
foreach ($r in $data) {
    $ws.Range("A$row").Value = $r[0]
    $ws.Range("B$row").Value = $r[1]
    $ws.Range("C$row").Value = $r[2]
    $ws.Range("D$row").Value = $r[3]
    $ws.Range("E$row").Value = $r[4]
    $ws.Range("F$row").Value = $r[5]
    $ws.Range("G$row").Value = "'" + $r[6]
    $ws.Range("G$row").Style = "Normal"
    $ws.Range("H$row").Value = $r[7]
    $ws.Range("J$row").Value = $r[8]
    $row++
}

# Final test row 18 - fully translated / reviewed sample row
$ws.Range("A18").Value = "TESTE-EAN"
$ws.Range("B18").Value = "TESTE"
$ws.Range("C18").Value = "MAQUILLAJE LABIOS"
$ws.Range("D18").Value = "Tiene PT"
$ws.Range("E18").Value = "Tiene ES"
$ws.Range("F18").Value = "Tiene IT"
$ws.Range("G18").Value = "'6"
$ws.Range("G18").Style = "Normal"
$ws.Range("H18").Value = "ML"
$ws.Range("I18").Value = "'"
$ws.Range("I18").Style = "Normal"
$ws.Range("J18").Value = "Revisado y Traducido"
$ws.Range("K18").Value = "TESTE Desc PT"
$ws.Range("L18").Value = "TESTE Uso PT"
$ws.Range("M18").Value = "TESTE prec PT"
$ws.Range("N18").Value = "TESTE +info PT"
$ws.Range("O18").Value = "TESTE desc IT"
$ws.Range("P18").Value = "TESTE uso IT"
$ws.Range("Q18").Value = "TESTE preca IT"
$ws.Range("R18").Value = "TESTE +info IT"
